$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 234, shifting existing rows 234:294 down to 235:295
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new record
$ws.Cells.Item(234, 1).Value = 11
$ws.Cells.Item(234, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(234, 3).Value = "Bíobío"
$ws.Cells.Item(234, 4).Value = 45093
$ws.Cells.Item(234, 5).Value = 8
$ws.Cells.Item(234, 6).Value = "Fruta"
$ws.Cells.Item(234, 7).Value = 100101
$ws.Cells.Item(234, 8).Value = "Berries"
$ws.Cells.Item(234, 9).Value = 100101007
$ws.Cells.Item(234, 10).Value = "Kiwi"
$ws.Cells.Item(234, 11).Value = "Hayward"
$ws.Cells.Item(234, 12).Value = "Primera"
$ws.Cells.Item(234, 13).Value = 190
$ws.Cells.Item(234, 14).Value = 8000
$ws.Cells.Item(234, 15).Value = 9000
$ws.Cells.Item(234, 16).Value = 8526
$ws.Cells.Item(234, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(234, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(234, 19).Value = 474
$ws.Cells.Item(234, 20).Value = 18

# Copy the date cell style (D column uses style index 2 for a date-like number format)
$ws.Range("D235").Copy()
$ws.Range("D234").PasteSpecial(-4122)  # xlPasteFormats
